# Update Digital Image Processing(DIP).pptx
#
# 1. Slide 1 (title slide): subtitle placeholder text "Subtitle" -> "SHILPA JOY"
# 2. Slide 4 ("Applications"): drop the trailing empty end-paragraph run
#    after "   Finger print recognition" (last bullet in the content list)

$p = $ppt.ActivePresentation

# --- 1. Title slide subtitle -------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "Subtitle") {
            $shape.TextFrame.TextRange.Text = "SHILPA JOY"
        }
    }
}

# --- 2. Applications slide bullet list --------------------------------------
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $shape = $slide4.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $paras = $tr.Paragraphs()
        $lastIndex = $paras.Count
        $lastPara = $tr.Paragraphs($lastIndex, 1)
        if ($lastPara.Text -match "Finger print recognition") {
            $lastPara.Text = "   Finger print recognition"
        }
    }
}
